$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the bold "header" style (previously only on A1:C1) onto the three
# new header cells so D1:F1 look like A1:C1 (bold, centered, bordered).
$ws.Range("C1").Copy($ws.Range("D1:F1"))

# --- Row 1: header labels (renamed + 3 new columns) ---
$ws.Range("A1").Value = "paciente"
$ws.Range("B1").Value = "tel.recado"
$ws.Range("C1").Value = "tel.celular"
$ws.Range("D1").Value = "message"
$ws.Range("E1").Value = "data.solicitacao"
$ws.Range("F1").Value = "diagnostico"

# --- Row 2 ---
$ws.Range("A2").Value = "João Silva"
$ws.Range("B2").Value = "11 - 9999 - 9999"
$ws.Range("C2").Value = "11 - 8888 - 8888"
$ws.Range("D2").Value = "Olá João! Lembrete de consulta."
$ws.Range("E2:F2").Value = ""
$ws.Range("E2:F2").Style = "Normal"

# --- Row 3 ---
$ws.Range("A3").Value = ""
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "11 - 9999 - 9999"
$ws.Range("C3").Value = "11 - 8888 - 8888"
$ws.Range("D3").Value = "Olá Lembrete de consulta."
$ws.Range("E3").Value = "29/08/2025"
$ws.Range("F3").Value = ""
$ws.Range("F3").Style = "Normal"

# --- Row 4 ---
$ws.Range("A4").Value = "Maria Santos"
$ws.Range("B4").Value = ""
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "11 - 7777 - 7777"
$ws.Range("D4").Value = "Olá Maria! Confirmação de horário."
$ws.Range("E4").Value = ""
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "cirurgia"

# --- Row 5 (new row) ---
$ws.Range("A5").Value = "Silva"
$ws.Range("B5:C5").Value = ""
$ws.Range("B5:C5").Style = "Normal"
$ws.Range("D5").Value = "Olá Silva! Lembrete de consulta."
$ws.Range("E5:F5").Value = ""
$ws.Range("E5:F5").Style = "Normal"
